# Applies cell updates to the cryptos price/volume table (cols D and E)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.932.18'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.619.58'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.499'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0618'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.843.59'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').Value = '1.623.44'
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.19%  '
$ws.Range('D16').Value = '25.931.08'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.84%  '
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('E33').Value = '  -2.52%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D36').Value = '1.124.68'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  -4.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  -1.24%  '
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('D42').Value = '1.755.07'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.08'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.10%  '
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('E51').Value = '  -0.47%  '
